$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "2020-06-30 00:00:00"
$ws.Range("I2").Value = 0.27
$ws.Range("J2").Value = 0.26
$ws.Range("K2").Value = 330780843.67
$ws.Range("L2").Value = 98051165.12
$ws.Range("M2").Value = 12.49
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = ""
$ws.Range("P2").Value = 1.388058541031
$ws.Range("Q2").Value = 0.256841734289
$ws.Range("R2").Value = 46.4484455978

$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "1"
$ws.Range("AC2").Value = "2020Q2"
$ws.Range("AD2").Value = "2020年 半年报"
$ws.Range("AE2").NumberFormat = "@"
$ws.Range("AE2").Value = "2020"
$ws.Range("AF2").Value = "半年报"
